$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C, rows 2 through 135 hold the "Förändrad" date serial.
# Update the date from serial 45202 (2023-10-03) to serial 45203 (2023-10-04)
# for every row in that range.
for ($row = 2; $row -le 135; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45202) {
        $cell.Value2 = 45203
    }
}
